$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---
# VALOR MORA total
$ws.Range("E11").Value = 180152
# Cant. Trabajadores
$ws.Range("C13").Value = 3
# Cant. Periodos
$ws.Range("F13").Value = 5

# --- Remove two rows from the worker detail table (22 -> 20 data rows stay 16-21) ---
# Delete from the bottom up so row indices of rows above stay stable.
$ws.Rows(22).Delete()
$ws.Rows(20).Delete()

# --- Rewrite the worker detail table contents (rows 16-21) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "30767353"
$ws.Range("D16").Value = "LENIS MARGARITA CRESPO MAYORAL"
$ws.Range("E16").Value = "1701"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 877803

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1003367663"
$ws.Range("D17").Value = "LUIS MIGUEL BALCEIRO LOPEZ"
$ws.Range("E17").Value = "1701"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "3829599"
$ws.Range("D18").Value = "CESAR AUGUSTO FERNANDEZ MARTINEZ"
$ws.Range("E18").Value = "1901"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "3829599"
$ws.Range("D19").Value = "CESAR AUGUSTO FERNANDEZ MARTINEZ"
$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "3829599"
$ws.Range("D20").Value = "CESAR AUGUSTO FERNANDEZ MARTINEZ"
$ws.Range("E20").Value = "1903"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "3829599"
$ws.Range("D21").Value = "CESAR AUGUSTO FERNANDEZ MARTINEZ"
$ws.Range("E21").Value = "1904"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 1423500
